$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns B through H (only column A remains)
$ws.Columns("B:H").Delete()

# The remaining column A keeps its original style; update its header text
# to what used to be column F's header ("input_Name")
$ws.Range("A1").Value = "input_Name"

# Row 2 in column A becomes an empty (text-typed) cell. Assigning a lone
# apostrophe stores it as a quote-prefixed empty string (no visible value);
# resetting the style afterwards drops the quote-prefix formatting so the
# cell is left as a plain empty text cell, matching the old (now removed)
# column F's empty cell.
$ws.Range("A2").Value = "'"
$ws.Range("A2").Style = "Normal"

# Column A width should match the old column F width (stored width 12)
$ws.Columns.Item(1).ColumnWidth = 11.17
